$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.752.76'
$ws.Range('E2').Value = '  -1.87%  '

# Row 3
$ws.Range('D3').Value = '1.897.28'
$ws.Range('E3').Value = '  -1.18%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9991'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.71%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.44%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9994'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.64%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4932'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.48%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3796'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.61%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07328'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.25%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9118'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.35%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.64'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.74%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07623'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.39%  '

# Row 13
$ws.Range('D13').Value = '1.892.61'
$ws.Range('E13').Value = '  -1.47%  '

# Row 14
$ws.Range('E14').Value = '  -1.68%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.659'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.20%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.18'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.21%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9997'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.61%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008733'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.96%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9987'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.70%  '

# Row 20
$ws.Range('D20').Value = '27.756.43'
$ws.Range('E20').Value = '  -1.89%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.84%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.124'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.01%  '

# Row 23
$ws.Range('D23').Value = '2.133.59'
$ws.Range('E23').Value = '  -1.19%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.76'
$ws.Range('D24').Style = 'Normal'

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.42%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.846'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.54%  '

# Row 27
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.174'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.80%  '

# Row 28
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.40'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.24%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '115.50'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.33%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.883'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.56%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08933'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.17%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.232'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.95%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.230'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.97%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7676'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.29%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.642'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.02%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.568'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.16%  '

# Row 37
$ws.Range('E37').Value = '  -0.77%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.100'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.98%  '

# Row 39
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05295'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.71%  '

# Row 40
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5500'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.84%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.990'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.41%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.909'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.42%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.581'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.43%  '

# Row 44
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1521'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.20%  '

# Row 45
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '112.29'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.94%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.62'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.42%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4806'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.86%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9992'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.68%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.637'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.63%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '67.50'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.93%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06051'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.67%  '
